# Updates sheet values from scheduled market-data runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2299.6
$ws.Range("I86").Value = 2374.5
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2374.5
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1251.5
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 2299.6
$ws.Range("I89").Value = 2374.5
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 11872.5
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -6256.5
$ws.Range("N89").Value = -21232

$ws.Range("H98").Value = 1088.6154
$ws.Range("J98").Value = 2000
$ws.Range("L98").Value = 2000
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 1088.6154
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 905.1429000000001
$ws.Range("I97").Value = 905.1429000000001
$ws.Range("K97").Value = 905.1429000000001
$ws.Range("M97").Value = -409.1429000000001

$ws.Range("H102").Value = 16442.217
$ws.Range("I102").Value = 19604.63
$ws.Range("J102").Value = 1420.75
$ws.Range("K102").Value = 19604.63
$ws.Range("L102").Value = 1420.75
$ws.Range("M102").Value = -17982.63
$ws.Range("N102").Value = -4664.75

$ws.Range("H122").Value = 4174.077
$ws.Range("I122").Value = 4141.3
$ws.Range("K122").Value = 12423.9
$ws.Range("M122").Value = -9973.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2503.923
$ws.Range("I80").Value = 660
$ws.Range("J80").Value = 3323.4443
$ws.Range("K80").Value = 660
$ws.Range("L80").Value = 3323.4443
$ws.Range("M80").Value = 338
$ws.Range("N80").Value = -5319.4443

$ws.Range("H83").Value = 2503.923
$ws.Range("I83").Value = 660
$ws.Range("J83").Value = 3323.4443
$ws.Range("K83").Value = 3300
$ws.Range("L83").Value = 16617.2215
$ws.Range("M83").Value = 1692
$ws.Range("N83").Value = -26601.2215

$ws.Range("H86").Value = 2103.4285
$ws.Range("I86").Value = 1928.3334
$ws.Range("J86").Value = 2418.6
$ws.Range("K86").Value = 1928.3334
$ws.Range("L86").Value = 2418.6
$ws.Range("M86").Value = -805.3334
$ws.Range("N86").Value = -4664.6

$ws.Range("H89").Value = 2103.4285
$ws.Range("I89").Value = 1928.3334
$ws.Range("J89").Value = 2418.6
$ws.Range("K89").Value = 9641.666999999999
$ws.Range("L89").Value = 12093
$ws.Range("M89").Value = -4025.666999999999
$ws.Range("N89").Value = -23325

$ws.Range("H99").Value = 20355.715
$ws.Range("I99").Value = 7079.8335
$ws.Range("J99").Value = 100011
$ws.Range("K99").Value = 7079.8335
$ws.Range("L99").Value = 100011
$ws.Range("M99").Value = -5581.8335
$ws.Range("N99").Value = -103007

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1010.25
$ws.Range("I22").Value = 1010.25
$ws.Range("K22").Value = 1010.25
$ws.Range("M22").Value = -660.25

$ws.Range("H41").Value = 19389
$ws.Range("I41").Value = 3059
$ws.Range("J41").Value = 24832.334
$ws.Range("K41").Value = 3059
$ws.Range("L41").Value = 24832.334
$ws.Range("M41").Value = -2631
$ws.Range("N41").Value = -25688.334

$ws.Range("H50").Value = 19500
$ws.Range("I50").Value = 19500
$ws.Range("K50").Value = 19500
$ws.Range("M50").Value = -18875

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -3855
$ws.Range("N59").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H68").Value = 94413.60000000001
$ws.Range("J68").Value = 94413.60000000001
$ws.Range("L68").Value = 94413.60000000001
$ws.Range("N68").Value = -95911.60000000001

$ws.Range("H71").Value = 94413.60000000001
$ws.Range("J71").Value = 94413.60000000001
$ws.Range("L71").Value = 283240.8
$ws.Range("N71").Value = -290728.8

$ws.Range("H74").Value = 49657
$ws.Range("J74").Value = 49657
$ws.Range("L74").Value = 49657
$ws.Range("N74").Value = -51405

$ws.Range("H77").Value = 49657
$ws.Range("J77").Value = 49657
$ws.Range("L77").Value = 148971
$ws.Range("N77").Value = -157707

$ws.Range("H99").Value = 2097.0833
$ws.Range("I99").Value = 1835.8
$ws.Range("J99").Value = 2283.7144
$ws.Range("K99").Value = 1835.8
$ws.Range("L99").Value = 2283.7144
$ws.Range("M99").Value = -337.8
$ws.Range("N99").Value = -5279.7144

$ws.Range("H122").Value = 3173.875
$ws.Range("J122").Value = 3878.6
$ws.Range("L122").Value = 11635.8
$ws.Range("N122").Value = -16535.8

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 2097.0833
$ws.Range("I126").Value = 1835.8
$ws.Range("J126").Value = 2283.7144
$ws.Range("K126").Value = 5507.4
$ws.Range("L126").Value = 6851.1432
$ws.Range("M126").Value = -3037.4
$ws.Range("N126").Value = -11791.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2949.25
$ws.Range("I20").Value = 40
$ws.Range("J20").Value = 3919
$ws.Range("K20").Value = 120
$ws.Range("L20").Value = 11757
$ws.Range("M20").Value = 107
$ws.Range("N20").Value = -12211

$ws.Range("H133").Value = 8619.352999999999
$ws.Range("I133").Value = 6868.5
$ws.Range("K133").Value = 20605.5
$ws.Range("M133").Value = -15545.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6036.625
$ws.Range("I102").Value = 5756.143
$ws.Range("K102").Value = 5756.143
$ws.Range("M102").Value = -4134.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2100.9092
$ws.Range("J46").Value = 2199.9
$ws.Range("L46").Value = 2199.9
$ws.Range("N46").Value = -2575.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 150678.4
$ws.Range("J119").Value = 150678.4
$ws.Range("L119").Value = 150678.4
$ws.Range("N119").Value = -160354.4

$ws.Range("H122").Value = 2886.0488
$ws.Range("I122").Value = 2478.743
$ws.Range("J122").Value = 5262
$ws.Range("K122").Value = 7436.228999999999
$ws.Range("L122").Value = 15786
$ws.Range("M122").Value = -4986.228999999999
$ws.Range("N122").Value = -20686
Write-Host "Applied market-data updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
